$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("enemyDatabase")

# Row 7 updates
$ws.Range("D7").Value = 100
$ws.Range("H7").Value = -0.2

# Row 10 updates
$ws.Range("D10").Value = 80
$ws.Range("I10").Value = 0.7
$ws.Range("J10").Value = 1.5
$ws.Range("O10").Value = 0
$ws.Range("T10").Value = "null"

# Update the active cell selection to O11
$ws.Activate()
$ws.Range("O11").Select()
